$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D (Price) cells keep their original text representation
# (values like "26.980.73" or "206.30" must remain literal text, not be
# reinterpreted as numbers by Excel).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.980.73"
$ws.Range("E2").Value = "  +0.41%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.552.75"
$ws.Range("E3").Value = "  -0.66%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.30"
$ws.Range("E5").Value = "  +0.34%  "
$ws.Range("E6").Value = "  -0.41%  "
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("E8").Value = "  +0.40%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "21.49"
$ws.Range("E9").Value = "  -1.45%  "
$ws.Range("E10").Value = "  -0.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0859"
$ws.Range("E11").Value = "  -0.60%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.771.55"
$ws.Range("E12").Value = "  -0.80%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.555.56"
$ws.Range("E13").Value = "  -0.56%  "
$ws.Range("E14").Value = "  -0.64%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.959.96"
$ws.Range("E16").Value = "  +0.30%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.63"
$ws.Range("E17").Value = "  +0.57%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "214.90"
$ws.Range("E18").Value = "  -0.21%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0686"
$ws.Range("E19").Value = "  +0.45%  "
$ws.Range("E20").Value = "  -2.10%  "
$ws.Range("E21").Value = "  +0.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.04"
$ws.Range("E22").Value = "  -2.35%  "
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("E24").Value = "  -3.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.56"
$ws.Range("E25").Value = "  -0.41%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.66"
$ws.Range("E26").Value = "  -0.41%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.86"
$ws.Range("E27").Value = "  -0.65%  "
$ws.Range("E28").Value = "  +0.23%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0461"
$ws.Range("E30").Value = "  -1.01%  "
$ws.Range("E31").Value = "  -1.07%  "
$ws.Range("E32").Value = "  +2.26%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.371.33"
$ws.Range("E33").Value = "  -1.64%  "
$ws.Range("E34").Value = "  +1.25%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.55"
$ws.Range("E35").Value = "  +1.64%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.966"
$ws.Range("E36").Value = "  +4.69%  "
$ws.Range("E37").Value = "  +0.24%  "
$ws.Range("E38").Value = "  +0.32%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.518"
$ws.Range("E39").Value = "  -2.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.808"
$ws.Range("E40").Value = "  -0.56%  "
$ws.Range("E41").Value = "  +0.29%  "
$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.987"
$ws.Range("E42").Value = "  -0.46%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.48"
$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.23"
$ws.Range("E44").Value = "  +1.91%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "63.75"
$ws.Range("E45").Value = "  +0.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.75"
$ws.Range("E46").Value = "  -1.74%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.685.21"
$ws.Range("E47").Value = "  -0.96%  "
$ws.Range("E48").Value = "  -2.95%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "86.27"
$ws.Range("E49").Value = "  -0.49%  "
$ws.Range("E50").Value = "  +0.85%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0954"
$ws.Range("E51").Value = "  +0.13%  "
